# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" values (column G) recomputed from the regenerated save_data.
$kValues = @{
    2  = 3
    3  = 1
    4  = 4
    5  = 6
    6  = 4
    7  = 4
    8  = 5
    9  = 4
    10 = 1
    11 = 6
    12 = 6
    13 = 5
    14 = 5
    15 = 3
    16 = 3
    17 = 4
    18 = 3
    19 = 2
    20 = 4
    21 = 2
    22 = 6
    23 = 5
    24 = 5
    25 = 3
    26 = 2
    27 = 2
    28 = 4
    29 = 4
    30 = 7
    31 = 3
    32 = 3
    33 = 2
    34 = 3
    36 = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
